# CapstoneHours.xlsx update
# Adds three new time-log entries ("Building my games Menu System") to
# Sheet1, covering 10/24/2014, 10/25/2014 and 10/26/2014 (Excel serials
# 41936-41938), for 2, 6 and 4 hours respectively. This mirrors the three
# existing rows already using that same task label/shared string and the
# same date-number-format style as the row directly above them (C63).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Sheet1" - the sheet touched by the edit

$label = "Building my games Menu System"

# Row, hours worked, and the date (as a raw Excel serial, same as the
# existing rows) for each new entry.
$newRows = @(
    @{ Row = 64; Hours = 2; DateSerial = 41936 },
    @{ Row = 65; Hours = 6; DateSerial = 41937 },
    @{ Row = 66; Hours = 4; DateSerial = 41938 }
)

# Grab the date cell style already used by row 63 (numFmt "date") so the
# new rows reuse the very same style instead of Excel inventing a new one.
$dateStyleSource = $ws.Cells.Item(63, 3)
$dateStyleSource.Copy()

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $label
    $ws.Cells.Item($r.Row, 2).Value = $r.Hours

    $dateCell = $ws.Cells.Item($r.Row, 3)
    $dateCell.PasteSpecial(-4122)   # xlPasteFormats: copy C63's number format/style
    $dateCell.Value = $r.DateSerial
}

# Make sure the Total Hours formula in F1 (SUM(B2,B3:B300)) picks up the
# newly added hours (138.5 -> 150.5).
$excel.Calculate()

# Match the final selection/scroll position left behind in the workbook.
$ws.Application.ActiveWindow.ScrollRow = 41
$ws.Range("C70").Select()
